$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): swap Result / SamplePortion in columns H and I ---
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# --- Row 2 (type markers): add unit to the float type used by H2/I2 ---
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"

# --- Row 3 (new): enum / description row for each column ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#Catalyseur"
